$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4719997
$ws.Range("J17").Value = 4719997
$ws.Range("L17").Value = 14159991
$ws.Range("N17").Value = -14160327

$ws.Range("H69").Value = 1466958.2
$ws.Range("I69").Value = 4600
$ws.Range("J69").Value = 1759429.9
$ws.Range("K69").Value = 13800
$ws.Range("L69").Value = 5278289.699999999
$ws.Range("M69").Value = -12926
$ws.Range("N69").Value = -5280037.699999999

$ws.Range("H72").Value = 1466958.2
$ws.Range("I72").Value = 4600
$ws.Range("J72").Value = 1759429.9
$ws.Range("K72").Value = 41400
$ws.Range("L72").Value = 15834869.1
$ws.Range("M72").Value = -37032
$ws.Range("N72").Value = -15843605.1

$ws.Range("H129").Value = 1063.9341
$ws.Range("I129").Value = 824.2857
$ws.Range("J129").Value = 1083.9048
$ws.Range("K129").Value = 2472.8571
$ws.Range("L129").Value = 3251.7144
$ws.Range("M129").Value = 2527.1429
$ws.Range("N129").Value = -13251.7144

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 644.2
$ws.Range("I2").Value = 509.30768
$ws.Range("K2").Value = 509.30768
$ws.Range("M2").Value = -396.30768

$ws.Range("H32").Value = 4659.26
$ws.Range("I32").Value = 3922.2197
$ws.Range("K32").Value = 3922.2197
$ws.Range("M32").Value = -3635.2197

$ws.Range("H61").Value = 237607.66
$ws.Range("I61").Value = 6596.5
$ws.Range("K61").Value = 6596.5
$ws.Range("M61").Value = -6384.5

$ws.Range("H63").Value = 76926170
$ws.Range("I63").Value = 111114410
$ws.Range("J63").Value = 2625
$ws.Range("K63").Value = 111114410
$ws.Range("L63").Value = 2625
$ws.Range("M63").Value = -111113724
$ws.Range("N63").Value = -3997

$ws.Range("H66").Value = 76926170
$ws.Range("I66").Value = 111114410
$ws.Range("J66").Value = 2625
$ws.Range("K66").Value = 555572050
$ws.Range("L66").Value = 13125
$ws.Range("M66").Value = -555568618
$ws.Range("N66").Value = -19989

$ws.Range("H102").Value = 6174532.5
$ws.Range("I102").Value = 9260762
$ws.Range("K102").Value = 9260762
$ws.Range("M102").Value = -9259140

$ws.Range("H116").Value = 644.2
$ws.Range("I116").Value = 509.30768
$ws.Range("K116").Value = 509.30768
$ws.Range("M116").Value = 1784.69232

$ws.Range("H136").Value = 237607.66
$ws.Range("I136").Value = 6596.5
$ws.Range("K136").Value = 19789.5
$ws.Range("M136").Value = -17239.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 644.2
$ws.Range("I3").Value = 509.30768
$ws.Range("K3").Value = 509.30768
$ws.Range("M3").Value = -395.30768

$ws.Range("H105").Value = 14206.842
$ws.Range("I105").Value = 36351.332
$ws.Range("J105").Value = 3986.3076
$ws.Range("K105").Value = 36351.332
$ws.Range("L105").Value = 3986.3076
$ws.Range("M105").Value = -34604.332
$ws.Range("N105").Value = -7480.3076

$ws.Range("H107").Value = 854.375
$ws.Range("I107").Value = 869.2857
$ws.Range("J107").Value = 750
$ws.Range("K107").Value = 869.2857
$ws.Range("L107").Value = 750
$ws.Range("M107").Value = 1050.7143
$ws.Range("N107").Value = -4590

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4552.3174
$ws.Range("I31").Value = 1873.317
$ws.Range("J31").Value = 9545
$ws.Range("K31").Value = 1873.317
$ws.Range("L31").Value = 9545
$ws.Range("M31").Value = -1578.317
$ws.Range("N31").Value = -10135

$ws.Range("H34").Value = 4552.3174
$ws.Range("I34").Value = 1873.317
$ws.Range("J34").Value = 9545
$ws.Range("K34").Value = 1873.317
$ws.Range("L34").Value = 9545
$ws.Range("M34").Value = -1671.317
$ws.Range("N34").Value = -9949

$ws.Range("H105").Value = 3232.3333
$ws.Range("I105").Value = 3201.45
$ws.Range("K105").Value = 3201.45
$ws.Range("M105").Value = -1454.45

$ws.Range("H118").Value = 39749.18
$ws.Range("J118").Value = 39749.18
$ws.Range("L118").Value = 39749.18
$ws.Range("N118").Value = -43063.18

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1244.6216
$ws.Range("I5").Value = 430.46155
$ws.Range("K5").Value = 1291.38465
$ws.Range("M5").Value = -1179.38465

$ws.Range("H86").Value = 893
$ws.Range("I86").Value = 893
$ws.Range("K86").Value = 2679
$ws.Range("M86").Value = -1493

$ws.Range("H89").Value = 893
$ws.Range("I89").Value = 893
$ws.Range("K89").Value = 8037
$ws.Range("M89").Value = -2109

$ws.Range("H131").Value = 1516057.1
$ws.Range("J131").Value = 1076.2554
$ws.Range("L131").Value = 3228.7662
$ws.Range("N131").Value = -13308.7662

$ws.Range("H135").Value = 1244.6216
$ws.Range("I135").Value = 430.46155
$ws.Range("K135").Value = 3874.15395
$ws.Range("M135").Value = -1339.15395

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5192.3706
$ws.Range("I70").Value = 5369.7
$ws.Range("J70").Value = 4685.7144
$ws.Range("K70").Value = 5369.7
$ws.Range("L70").Value = 4685.7144
$ws.Range("M70").Value = -5099.7
$ws.Range("N70").Value = -5225.7144

$ws.Range("H73").Value = 5192.3706
$ws.Range("I73").Value = 5369.7
$ws.Range("J73").Value = 4685.7144
$ws.Range("K73").Value = 5369.7
$ws.Range("L73").Value = 4685.7144
$ws.Range("M73").Value = -4433.7
$ws.Range("N73").Value = -6557.7144

$ws.Range("H107").Value = 871.4761999999999
$ws.Range("I107").Value = 212.28572
$ws.Range("J107").Value = 2189.8572
$ws.Range("K107").Value = 212.28572
$ws.Range("L107").Value = 2189.8572
$ws.Range("M107").Value = 1707.71428
$ws.Range("N107").Value = -6029.8572

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1712.5151
$ws.Range("I22").Value = 401.625
$ws.Range("J22").Value = 2132
$ws.Range("K22").Value = 401.625
$ws.Range("L22").Value = 2132
$ws.Range("M22").Value = -106.625
$ws.Range("N22").Value = -2722

$ws.Range("H27").Value = 1712.5151
$ws.Range("I27").Value = 401.625
$ws.Range("J27").Value = 2132
$ws.Range("K27").Value = 401.625
$ws.Range("L27").Value = 2132
$ws.Range("M27").Value = -294.625
$ws.Range("N27").Value = -2346

$ws.Range("H68").Value = 90911160
$ws.Range("I68").Value = 2250.25
$ws.Range("J68").Value = 333334940
$ws.Range("K68").Value = 2250.25
$ws.Range("L68").Value = 333334940
$ws.Range("M68").Value = -1501.25
$ws.Range("N68").Value = -333336438

$ws.Range("H71").Value = 90911160
$ws.Range("I71").Value = 2250.25
$ws.Range("J71").Value = 333334940
$ws.Range("K71").Value = 11251.25
$ws.Range("L71").Value = 1666674700
$ws.Range("M71").Value = -7507.25
$ws.Range("N71").Value = -1666682188

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1623.3
$ws.Range("I132").Value = 1081.742
$ws.Range("J132").Value = 3488.6667
$ws.Range("K132").Value = 3245.226
$ws.Range("L132").Value = 10466.0001
$ws.Range("M132").Value = -715.2259999999997
$ws.Range("N132").Value = -15526.0001
